$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.38"
$ws.Range("E2").Value = "'-0.03%"
$ws.Range("D3").Value = "'41.07"
$ws.Range("E3").Value = "'0.40%"
$ws.Range("D4").Value = "'5.201"
$ws.Range("E4").Value = "'1.70%"
$ws.Range("D5").Value = "'0.07673"
$ws.Range("E5").Value = "'0.63%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.295"
$ws.Range("E6").Value = "'1.21%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.633"
$ws.Range("E7").Value = "'1.64%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9148"
$ws.Range("E8").Value = "'1.29%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.430"
$ws.Range("E9").Value = "'0.16%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1229"
$ws.Range("E10").Value = "'10.59%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1824"
$ws.Range("E11").Value = "'2.56%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09160"
$ws.Range("E12").Value = "'1.12%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04263"
$ws.Range("E13").Value = "'1.40%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1051"
$ws.Range("E14").Value = "'-0.05%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001258"
$ws.Range("E15").Value = "'0.09%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005852"
$ws.Range("E16").Value = "'1.85%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007509"
$ws.Range("E17").Value = "'1,903.55%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.343"
$ws.Range("E18").Value = "'-0.23%"
$ws.Range("D20").Value = "'7.348"
$ws.Range("E20").Value = "'11.77%"
$ws.Range("E21").Value = "'1.24%"
$ws.Range("D22").Value = "'0.2712"
$ws.Range("E22").Value = "'-4.17%"
$ws.Range("D23").Value = "'0.04027"
$ws.Range("E23").Value = "'-0.90%"
$ws.Range("E24").Value = "'2.74%"
$ws.Range("E25").Value = "'6.59%"
$ws.Range("E26").Value = "'0.08%"
$ws.Range("D38").Value = "'0.02497"
$ws.Range("E38").Value = "'3.65%"
$ws.Range("D39").Value = "'0.05316"
$ws.Range("E39").Value = "'2.74%"
$ws.Range("D40").Value = "'0.007840"
$ws.Range("E40").Value = "'1.05%"
$ws.Range("E41").Value = "'0.99%"
$ws.Range("D42").Value = "'0.006662"
$ws.Range("E42").Value = "'-5.47%"
$ws.Range("D43").Value = "'0.001861"
$ws.Range("E43").Value = "'-4.54%"
$ws.Range("D44").Value = "'0.008010"
$ws.Range("E44").Value = "'-8.94%"
$ws.Range("D45").Value = "'0.3066"
$ws.Range("E45").Value = "'-0.43%"
$ws.Range("D46").Value = "'0.00006741"
$ws.Range("E46").Value = "'-3.05%"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("D48").Value = "'0.2986"
$ws.Range("E48").Value = "'868.07%"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("E51").Value = "'0.09%"
